$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 credential values (newly added iAuthor TC's)
$ws.Range("A2").Value = "NFspP583"
$ws.Range("B2").Value = 23110948
$ws.Range("C2").Value = "hblzzlx44"
$ws.Range("D2").Value = "qj8#&D4M"
$ws.Range("F2").Value = "HCYYylXe"
$ws.Range("G2").Value = "NfWh"
